$d = $word.ActiveDocument

# --- 1) "En cas de problemes graves ... son CdP." : split "CdP" into its own run ---
$r = $d.Content
$f = $r.Find
$f.Execute("CdP")
$cdpStart = $r.Start
$cdpEnd = $r.End
$cdpRange = $d.Range($cdpStart, $cdpEnd)
# Forces a run boundary around "CdP" without altering the resolved formatting
# (same trick Word uses internally when it (re)splits runs on save).
$cdpRange.Bold = 1
$cdpRange.Bold = 0

# --- 2) "Mostly Human" -> "Shadow scan" ---
$d.Content.Find.Execute("Mostly Human", $true, $false, $false, $false, $false, $true, 1, $false, "Shadow scan", 2)

# --- 3) "Le prof peut modifier la bibliothèque" -> "Management des ressource bloquées" ---
$d.Content.Find.Execute("Le prof peut modifier la bibliothèque", $true, $false, $false, $false, $false, $true, 1, $false, "Management des ressource bloquées", 2)

# --- 4) "Le prof peut créer une sous liste" -> "Robustesse de la connexion entre les machines" ---
$d.Content.Find.Execute("Le prof peut créer une sous liste", $true, $false, $false, $false, $false, $true, 1, $false, "Robustesse de la connexion entre les machines", 2)
